# Update and revise cost inputs
# - Add NNS
# - Add input costs for decentralisation
# - Change start-up cost from annual to total

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# program_number_tests_per_tb_presentation: value revised 2 -> 50
$ws.Range("B45").Value = 50

# Remove stray review notes in column F next to the NNS rows
$ws.Range("F43").ClearContents()
$ws.Range("F44").ClearContents()

# econ_startupcost_ipt: annual -> total start-up cost
$ws.Range("B119").Value = 70800

# econ_startupcost_xpert: annual -> total start-up cost
$ws.Range("B124").Value = 311038.5

# econ_startupcost_treatment_support: annual -> total start-up cost
$ws.Range("B129").Value = 10407

# Drop the "Average cost of 3 years" notes next to the smearacf / xpertacf inflection costs
$ws.Range("E133").ClearContents()

# econ_startupcost_smearacf: annual -> total start-up cost
$ws.Range("B134").Value = 277254.63

$ws.Range("E138").ClearContents()

# econ_startupcost_xpertacf: annual -> total start-up cost
$ws.Range("B139").Value = 567148.2

# Input costs for decentralisation (was using dummy/placeholder values)
$ws.Range("B142").Value = 100
$ws.Range("E142").ClearContents()

$ws.Range("E143").ClearContents()

# econ_startupcost_decentralisation: annual -> total start-up cost
$ws.Range("B144").Value = 1000000
$ws.Range("E145").ClearContents()

# econ_saturation_decentralisation revised
$ws.Range("B146").Value = 0.8
$ws.Range("E146").ClearContents()

# Update the visible selection to reflect where the author left off editing
$ws.Activate()
$ws.Range("B145").Select()
